# Correction Planning et visu Appel Micro
# Update the "Données" sheet: shorten the event-id labels to their short
# codes and correct the "Temps/Series" and "Battement Epreuve" values.
# The "Visuel_Planning" sheet reads these via formulas and recalculates
# automatically.

$wb  = $excel.ActiveWorkbook
$dat = $wb.Worksheets.Item("Données")
$viz = $wb.Worksheets.Item("Visuel_Planning")

# --- Id Epreuve column (A) : long labels -> short codes -------------------
$dat.Range("A2").Value = "sta"
$dat.Range("A3").Value = "dwf"
$dat.Range("A4").Value = "spd"
$dat.Range("A5").Value = "dnf"
$dat.Range("A6").Value = 1650

# --- Temps/Series column (E) -----------------------------------------------
$dat.Range("E2").Value = 12
$dat.Range("E3").Value = 9
$dat.Range("E4").Value = 7
$dat.Range("E5").Value = 9
$dat.Range("E6").Value = 23

# --- Battement Epreuve column (G) ------------------------------------------
$dat.Range("G3").Value = 39
$dat.Range("G5").Value = 9

# --- Fix the off-by-one label references on "Visuel_Planning" --------------
# Each block of 4/3 rows shows the label for one "Données" row; the formulas
# for blocks 2-5 incorrectly all pointed at the row used by block 1/2. Point
# every block at its own "Données" row.
$viz.Range("B5").Formula  = "=Données!A3"
$viz.Range("B9").Formula  = "=Données!A4"
$viz.Range("B13").Formula = "=Données!A5"
$viz.Range("B17").Formula = "=Données!A6"

# --- Restore the cell selections as left by the author ---------------------
$dat.Range("G6").Select()
$viz.Range("E4").Select()
$viz.Activate()
